# proteiNorm_IDeA.pptx edit
# Commit: "updated figures, reduced data / using top3 method now / only providing batch 3 and 4 now"
#
# The canonical diff shows slide 3 ("Filter outliers") being removed from the
# deck entirely; every other visible change in the raw XML diff (cached
# datetimeFigureOut / slidenum field text, relationship-id renumbering, etc.)
# is a mechanical side effect of PowerPoint re-saving the file after that
# slide is deleted, not a distinct edit.

$p = $ppt.ActivePresentation

# Slide 3 is "Filter outliers" in the original deck order.
$s = $p.Slides.Item(3)
$s.Delete()
